$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new weather-log rows (12-24). Temp/humidity columns (C/D) hold
# numeric-looking text in the source data, so force NumberFormat to "@"
# (text) before assigning the value to keep them stored as strings
# instead of being auto-coerced to numbers.

$ws.Range("A12").Value = "São Paulo"
$ws.Range("B12").Value = "26/08/2025 21:49"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "16"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "90"
$ws.Range("E12").Value = "Nublado"

$ws.Range("A13").Value = "São Paulo"
$ws.Range("B13").Value = "26/08/2025 22:22"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "16"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "90"
$ws.Range("E13").Value = "Nublado"

$ws.Range("A14").Value = "São Paulo"
$ws.Range("B14").Value = "26/08/2025 23:22"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "16"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "91"
$ws.Range("E14").Value = "Nublado"

$ws.Range("A15").Value = "São Paulo"
$ws.Range("B15").Value = "26/08/2025 23:25"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "16"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92"
$ws.Range("E15").Value = "Nublado"

$ws.Range("A16").Value = "São Paulo"
$ws.Range("B16").Value = "26/08/2025 23:27"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "16"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92"
$ws.Range("E16").Value = "Nublado"

$ws.Range("A17").Value = "São Paulo"
$ws.Range("B17").Value = "26/08/2025 23:58"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "16"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94"
$ws.Range("E17").Value = "Nublado"

$ws.Range("A18").Value = "São Paulo"
$ws.Range("B18").Value = "27/08/2025 00:02"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "16"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94"
$ws.Range("E18").Value = "Nublado"

$ws.Range("A19").Value = "BRASILIA"
$ws.Range("B19").Value = "27/08/2025 00:03"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "19"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "54"
$ws.Range("E19").Value = "Céu limpo com poucas nuvens"

$ws.Range("A20").Value = "BARSILIA"
$ws.Range("B20").Value = "27/08/2025 00:05"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "19"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "54"
$ws.Range("E20").Value = "Céu limpo com poucas nuvens"

$ws.Range("A21").Value = "BARSILIA"
$ws.Range("B21").Value = "27/08/2025 00:06"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "19"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "54"
$ws.Range("E21").Value = "Céu limpo com poucas nuvens"

$ws.Range("A22").Value = "BARSILIA"
$ws.Range("B22").Value = "27/08/2025 00:07"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "19"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "54"
$ws.Range("E22").Value = "Céu limpo com poucas nuvens"

$ws.Range("A23").Value = "São Paulo"
$ws.Range("B23").Value = "27/08/2025 00:10"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "16"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94"
$ws.Range("E23").Value = "Nublado"

$ws.Range("A24").Value = "São Paulo"
$ws.Range("B24").Value = "27/08/2025 00:25"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "16"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "94"
$ws.Range("E24").Value = "Nublado"

# Widen the Status column (E) from 22 to 29 characters. The engine's
# ColumnWidth setter round-trips through a pixel conversion that adds a
# fixed 5/6-character padding on save, so subtract it here to land on an
# exact stored width of 29.
$ws.Columns.Item(5).ColumnWidth = 29 - 5/6
